# Journal Entries - add a new entry row (row 8): date 2023-11-05 with
# blank "title"/"entry" fields, mirroring the existing table layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to stay plain text so the date-looking string isn't
# auto-converted into a date serial number, then drop the temporary
# number-format tweak so the cell keeps the sheet's default style.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2023-11-05"
$ws.Range("A8").ClearFormats()

# New entry has no title/body yet.
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
